# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-name suffixes used in row 1 to the
#   new version tags "_FV2410" / "_FV2504".
# - Turn the data range A1:U65 into a proper Excel Table ("Table1").
# - Freeze the header row (row 1) in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rewrite the 21 header cells on row 1 (A1:U1).
$headers = @(
    "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410", "Segment ID_FV2410",
    "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410", "Bedingungsausdruck_FV2410", "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504", "Segmentgruppe_FV2504", "Segment_FV2504", "Datenelement_FV2504", "Segment ID_FV2504",
    "Code_FV2504", "Qualifier_FV2504", "Beschreibung_FV2504", "Bedingungsausdruck_FV2504", "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# 2) Convert A1:U65 into an Excel Table named "Table1" (adds xl/tables/table1.xml,
#    the worksheet <tableParts> entry and the part relationship).
$range = $ws.Range("A1:U65")
$lo = $ws.ListObjects.Add(1, $range, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# 3) Freeze panes above row 2 (keeps the header row visible while scrolling).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
